$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format first so numeric-looking strings
# (e.g. "305.34", "0.0810") are preserved exactly as text, matching
# the inlineStr cell type used in the source workbook.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = '44.235.39'
$ws.Range("E2").Value = '  +0.90%  '

$ws.Range("D3").Value = '2.241.39'
$ws.Range("E3").Value = '  +0.26%  '

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").Value = '305.34'
$ws.Range("E5").Value = '  -2.93%  '

$ws.Range("D6").Value = '95.44'
$ws.Range("E6").Value = '  -3.01%  '

$ws.Range("D7").Value = '0.571'
$ws.Range("E7").Value = '  -0.17%  '

$ws.Range("E8").Value = '  +0.20%  '

$ws.Range("D9").Value = '0.525'
$ws.Range("E9").Value = '  -1.30%  '

$ws.Range("D10").Value = '34.76'
$ws.Range("E10").Value = '  -2.90%  '

$ws.Range("D11").Value = '0.0810'
$ws.Range("E11").Value = '  -1.28%  '

$ws.Range("D12").Value = '7.22'
$ws.Range("E12").Value = '  -1.72%  '

$ws.Range("E13").Value = '  -0.11%  '

$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '2.583.07'
$ws.Range("E14").Value = '  +0.21%  '

$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '2.318.84'
$ws.Range("E15").Value = '  +3.64%  '

$ws.Range("D16").Value = '0.829'
$ws.Range("E16").Value = '  -0.98%  '

$ws.Range("D17").Value = '13.53'
$ws.Range("E17").Value = '  -2.89%  '

$ws.Range("D18").Value = '44.010.68'
$ws.Range("E18").Value = '  +0.66%  '

$ws.Range("D19").Value = '0.0₃0963'
$ws.Range("E19").Value = '  -0.04%  '

$ws.Range("D20").Value = '6.35'
$ws.Range("E20").Value = '  +1.13%  '

$ws.Range("D21").Value = '12.03'
$ws.Range("E21").Value = '  -8.12%  '

$ws.Range("D22").Value = '65.52'
$ws.Range("E22").Value = '  -0.57%  '

$ws.Range("B23").Value = 'PancakeSwap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D23").Value = '3.12'
$ws.Range("E23").Value = '  +4.53%  '

$ws.Range("B24").Value = 'BitcoinCash'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D24").Value = '237.78'
$ws.Range("E24").Value = '  +0.65%  '

$ws.Range("D25").Value = '2.00'
$ws.Range("E25").Value = '  -1.14%  '

$ws.Range("E26").Value = '  -0.02%  '

$ws.Range("D27").Value = '38.28'
$ws.Range("E27").Value = '  +4.72%  '

$ws.Range("D28").Value = '9.90'
$ws.Range("E28").Value = '  -1.66%  '

$ws.Range("D29").Value = '2.17'
$ws.Range("E29").Value = '  +0.75%  '

$ws.Range("D30").Value = '20.02'
$ws.Range("E30").Value = '  +0.35%  '

$ws.Range("D31").Value = '5.85'
$ws.Range("E31").Value = '  -1.81%  '

$ws.Range("D32").Value = '152.31'
$ws.Range("E32").Value = '  -2.64%  '

$ws.Range("D33").Value = '0.0791'
$ws.Range("E33").Value = '  -4.75%  '

$ws.Range("B34").Value = 'WEMIXToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D34").Value = '2.61'
$ws.Range("E34").Value = '  -1.04%  '

$ws.Range("B35").Value = 'LidoDAOToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D35").Value = '3.22'
$ws.Range("E35").Value = '  -3.14%  '

$ws.Range("E36").Value = '  +1.82%  '

$ws.Range("D37").Value = '0.106'
$ws.Range("E37").Value = '  -2.41%  '

$ws.Range("E38").Value = '  -7.95%  '

$ws.Range("D39").Value = '3.56'
$ws.Range("E39").Value = '  +1.20%  '

$ws.Range("D40").Value = '3.84'
$ws.Range("E40").Value = '  -3.91%  '

$ws.Range("D41").Value = '14.41'
$ws.Range("E41").Value = '  -6.88%  '

$ws.Range("D42").Value = '0.0298'
$ws.Range("E42").Value = '  -2.68%  '

$ws.Range("E43").Value = '  +0.25%  '

$ws.Range("D44").Value = '1.752.82'
$ws.Range("E44").Value = '  +2.85%  '

$ws.Range("D45").Value = '82.56'
$ws.Range("E45").Value = '  +0.36%  '

$ws.Range("D46").Value = '0.190'
$ws.Range("E46").Value = '  -1.81%  '

$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").Value = '99.83'
$ws.Range("E47").Value = '  -1.70%  '

$ws.Range("B48").Value = 'THORChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D48").Value = '4.94'
$ws.Range("E48").Value = '  -4.05%  '

$ws.Range("D49").Value = '1.58'
$ws.Range("E49").Value = '  -1.68%  '

$ws.Range("D50").Value = '8.10'
$ws.Range("E50").Value = '  -0.50%  '

$ws.Range("D51").Value = '54.56'
$ws.Range("E51").Value = '  -2.89%  '

# Restore default style on column D so no stray formatting remains
$dRange.Style = "Normal"
